$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 8).Value = 1218
$ws.Cells.Item(2, 10).Value = 3011
$ws.Cells.Item(2, 12).Value = 3011
$ws.Cells.Item(2, 14).Value = -3237
$ws.Cells.Item(9, 8).Value = 284.2
$ws.Cells.Item(9, 9).Value = 330.25
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 330.25
$ws.Cells.Item(9, 12).Value = 100
$ws.Cells.Item(9, 13).Value = -161.25
$ws.Cells.Item(9, 14).Value = -438
$ws.Cells.Item(31, 8).Value = 197.14285
$ws.Cells.Item(31, 10).Value = 1200
$ws.Cells.Item(31, 12).Value = 3600
$ws.Cells.Item(31, 14).Value = -4060
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(55, 14).ClearContents()
$ws.Cells.Item(74, 8).Value = 7975.2856
$ws.Cells.Item(74, 9).Value = 7456.75
$ws.Cells.Item(74, 11).Value = 7456.75
$ws.Cells.Item(74, 13).Value = -6520.75
$ws.Cells.Item(77, 8).Value = 7975.2856
$ws.Cells.Item(77, 9).Value = 7456.75
$ws.Cells.Item(77, 11).Value = 37283.75
$ws.Cells.Item(77, 13).Value = -32603.75
$ws.Cells.Item(132, 8).Value = 13257.177
$ws.Cells.Item(132, 9).Value = 14624.8
$ws.Cells.Item(132, 11).Value = 43874.39999999999
$ws.Cells.Item(132, 13).Value = -41344.39999999999
$ws.Cells.Item(135, 8).Value = 2638.4
$ws.Cells.Item(135, 9).Value = 2464.75
$ws.Cells.Item(135, 11).Value = 22182.75
$ws.Cells.Item(135, 13).Value = -19647.75

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(32, 8).Value = 6852.5
$ws.Cells.Item(32, 9).Value = 3869.5833
$ws.Cells.Item(32, 10).Value = 24750
$ws.Cells.Item(32, 11).Value = 3869.5833
$ws.Cells.Item(32, 12).Value = 24750
$ws.Cells.Item(32, 13).Value = -3582.5833
$ws.Cells.Item(32, 14).Value = -25324
$ws.Cells.Item(44, 8).Value = 11614.35
$ws.Cells.Item(44, 10).Value = 11614.35
$ws.Cells.Item(44, 12).Value = 11614.35
$ws.Cells.Item(44, 14).Value = -12590.35
$ws.Cells.Item(55, 8).Value = 58999.5
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 58999.5
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 12).Value = 58999.5
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(55, 14).Value = -59629.5
$ws.Cells.Item(74, 8).Value = 8463.333000000001
$ws.Cells.Item(74, 9).Value = 8419.777
$ws.Cells.Item(74, 10).Value = 8594
$ws.Cells.Item(74, 11).Value = 8419.777
$ws.Cells.Item(74, 12).Value = 8594
$ws.Cells.Item(74, 13).Value = -7545.777
$ws.Cells.Item(74, 14).Value = -10342
$ws.Cells.Item(77, 8).Value = 8463.333000000001
$ws.Cells.Item(77, 9).Value = 8419.777
$ws.Cells.Item(77, 10).Value = 8594
$ws.Cells.Item(77, 11).Value = 42098.885
$ws.Cells.Item(77, 12).Value = 42970
$ws.Cells.Item(77, 13).Value = -37730.885
$ws.Cells.Item(77, 14).Value = -51706
$ws.Cells.Item(98, 8).Value = 10000
$ws.Cells.Item(98, 10).Value = 10000
$ws.Cells.Item(98, 12).Value = 10000
$ws.Cells.Item(98, 14).Value = -15990
$ws.Cells.Item(122, 8).Value = 799.5
$ws.Cells.Item(122, 9).Value = 799.5
$ws.Cells.Item(122, 11).Value = 2398.5
$ws.Cells.Item(122, 13).Value = 51.5
$ws.Cells.Item(132, 8).Value = 2135.2727
$ws.Cells.Item(132, 9).Value = 2248.9
$ws.Cells.Item(132, 10).Value = 999
$ws.Cells.Item(132, 11).Value = 6746.700000000001
$ws.Cells.Item(132, 12).Value = 2997
$ws.Cells.Item(132, 13).Value = -4216.700000000001
$ws.Cells.Item(132, 14).Value = -8057

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(22, 8).Value = 184
$ws.Cells.Item(22, 9).Value = 170.8
$ws.Cells.Item(22, 11).Value = 170.8
$ws.Cells.Item(22, 13).Value = 2.199999999999989
$ws.Cells.Item(94, 8).Value = 260.5
$ws.Cells.Item(94, 9).Value = 237.5
$ws.Cells.Item(94, 11).Value = 237.5
$ws.Cells.Item(94, 13).Value = 213.5

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(7, 8).Value = 67.61539
$ws.Cells.Item(7, 9).Value = 45
$ws.Cells.Item(7, 10).Value = 143
$ws.Cells.Item(7, 11).Value = 45
$ws.Cells.Item(7, 12).Value = 143
$ws.Cells.Item(7, 13).Value = 68
$ws.Cells.Item(7, 14).Value = -369
$ws.Cells.Item(22, 8).Value = 1635.25
$ws.Cells.Item(22, 9).Value = 1019.2222
$ws.Cells.Item(22, 10).Value = 3483.3333
$ws.Cells.Item(22, 11).Value = 1019.2222
$ws.Cells.Item(22, 12).Value = 3483.3333
$ws.Cells.Item(22, 13).Value = -669.2222
$ws.Cells.Item(22, 14).Value = -4183.3333
$ws.Cells.Item(43, 8).Value = 13008
$ws.Cells.Item(43, 10).Value = 13008
$ws.Cells.Item(43, 12).Value = 13008
$ws.Cells.Item(43, 14).Value = -13376
$ws.Cells.Item(99, 8).Value = 2534.2942
$ws.Cells.Item(99, 10).Value = 2643.1667
$ws.Cells.Item(99, 12).Value = 2643.1667
$ws.Cells.Item(99, 14).Value = -5639.1667
$ws.Cells.Item(101, 8).Value = 13008
$ws.Cells.Item(101, 10).Value = 13008
$ws.Cells.Item(101, 12).Value = 13008
$ws.Cells.Item(101, 14).Value = -19498
$ws.Cells.Item(126, 8).Value = 2534.2942
$ws.Cells.Item(126, 10).Value = 2643.1667
$ws.Cells.Item(126, 12).Value = 7929.500100000001
$ws.Cells.Item(126, 14).Value = -12869.5001
$ws.Cells.Item(132, 8).Value = 4399.222
$ws.Cells.Item(132, 9).Value = 4008.3635
$ws.Cells.Item(132, 10).Value = 5013.4287
$ws.Cells.Item(132, 11).Value = 12025.0905
$ws.Cells.Item(132, 12).Value = 15040.2861
$ws.Cells.Item(132, 13).Value = -9495.0905
$ws.Cells.Item(132, 14).Value = -20100.2861

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(55, 8).Value = 6140.067
$ws.Cells.Item(55, 10).Value = 6845.923
$ws.Cells.Item(55, 12).Value = 20537.769
$ws.Cells.Item(55, 14).Value = -20891.769
$ws.Cells.Item(122, 8).Value = 1074.5
$ws.Cells.Item(122, 9).Value = 1346.5
$ws.Cells.Item(122, 10).Value = 802.5
$ws.Cells.Item(122, 11).Value = 12118.5
$ws.Cells.Item(122, 12).Value = 7222.5
$ws.Cells.Item(122, 13).Value = -9668.5
$ws.Cells.Item(122, 14).Value = -12122.5
$ws.Cells.Item(131, 8).Value = 2466
$ws.Cells.Item(131, 10).Value = 2466
$ws.Cells.Item(131, 12).Value = 7398
$ws.Cells.Item(131, 14).Value = -17478

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(63, 8).Value = 49990
$ws.Cells.Item(63, 9).Value = 49990
$ws.Cells.Item(63, 11).Value = 49990
$ws.Cells.Item(63, 13).Value = -49304
$ws.Cells.Item(66, 8).Value = 49990
$ws.Cells.Item(66, 9).Value = 49990
$ws.Cells.Item(66, 11).Value = 149970
$ws.Cells.Item(66, 13).Value = -146538
$ws.Cells.Item(113, 8).Value = 10001.5
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 10001.5
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 10001.5
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).Value = -14341.5
$ws.Cells.Item(122, 8).Value = 6123
$ws.Cells.Item(122, 9).Value = 5187.5
$ws.Cells.Item(122, 10).Value = 7994
$ws.Cells.Item(122, 11).Value = 15562.5
$ws.Cells.Item(122, 12).Value = 23982
$ws.Cells.Item(122, 13).Value = -13112.5
$ws.Cells.Item(122, 14).Value = -28882

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(22, 8).Value = 903.4545000000001
$ws.Cells.Item(22, 9).Value = 948.6667
$ws.Cells.Item(22, 11).Value = 948.6667
$ws.Cells.Item(22, 13).Value = -653.6667
$ws.Cells.Item(27, 8).Value = 903.4545000000001
$ws.Cells.Item(27, 9).Value = 948.6667
$ws.Cells.Item(27, 11).Value = 948.6667
$ws.Cells.Item(27, 13).Value = -841.6667
$ws.Cells.Item(42, 8).Value = 39749.5
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 39749.5
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 39749.5
$ws.Cells.Item(42, 13).ClearContents()
$ws.Cells.Item(42, 14).Value = -40875.5
$ws.Cells.Item(46, 8).Value = 4867.619
$ws.Cells.Item(46, 9).Value = 4037.5
$ws.Cells.Item(46, 11).Value = 4037.5
$ws.Cells.Item(46, 13).Value = -3849.5
$ws.Cells.Item(47, 8).Value = 10001
$ws.Cells.Item(47, 9).Value = 10001
$ws.Cells.Item(47, 11).Value = 10001
$ws.Cells.Item(47, 13).Value = -9511
$ws.Cells.Item(49, 8).Value = 39749.5
$ws.Cells.Item(49, 9).Value = 0
$ws.Cells.Item(49, 10).Value = 39749.5
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(49, 12).Value = 39749.5
$ws.Cells.Item(49, 13).ClearContents()
$ws.Cells.Item(49, 14).Value = -40043.5
$ws.Cells.Item(52, 8).Value = 10001
$ws.Cells.Item(52, 9).Value = 10001
$ws.Cells.Item(52, 11).Value = 10001
$ws.Cells.Item(52, 13).Value = -9768
$ws.Cells.Item(68, 8).Value = 7399.778
$ws.Cells.Item(68, 9).Value = 4399.5
$ws.Cells.Item(68, 11).Value = 4399.5
$ws.Cells.Item(68, 13).Value = -3650.5
$ws.Cells.Item(71, 8).Value = 7399.778
$ws.Cells.Item(71, 9).Value = 4399.5
$ws.Cells.Item(71, 11).Value = 21997.5
$ws.Cells.Item(71, 13).Value = -18253.5
$ws.Cells.Item(93, 8).Value = 1206.6666
$ws.Cells.Item(93, 9).Value = 1206.6666
$ws.Cells.Item(93, 11).Value = 1206.6666
$ws.Cells.Item(93, 13).Value = 41.33339999999998
$ws.Cells.Item(101, 8).Value = 26749.5
$ws.Cells.Item(101, 10).Value = 26749.5
$ws.Cells.Item(101, 12).Value = 26749.5
$ws.Cells.Item(101, 14).Value = -33239.5
$ws.Cells.Item(105, 8).Value = 32374.25
$ws.Cells.Item(105, 10).Value = 32374.25
$ws.Cells.Item(105, 12).Value = 32374.25
$ws.Cells.Item(105, 14).Value = -39362.25

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(52, 8).Value = 38367.285
$ws.Cells.Item(52, 10).Value = 57124.25
$ws.Cells.Item(52, 12).Value = 57124.25
$ws.Cells.Item(52, 14).Value = -57576.25
